$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header typo: "Assets Requied" -> "Assets Required" ---
$ws.Range("F2").Value = "Assets Required"

# --- Row 7 (LandOnEnemy): replace Assets Required text and add Status / Notes ---
$ws.Range("F7").Value = "Jawharp multi instrument, bubbling mud sample"
$ws.Range("H7").Value = "5 different possible jawharp twangs, -/+ .25 semitones"
$ws.Range("G7").Value = "Iteration 1"

# --- Widen Status column (G) ---
$ws.Columns.Item(7).ColumnWidth = 13.14

# --- Update current selection in the sheet view ---
$ws.Range("F16").Select()
